# "Generate Report for Handoff"
# The 631a7a47-57af-43cf-bb8c-79ccb0c4fd9e file has completed handoff: its
# status flips from "In Translation" to "Ready for handoff" on every sheet
# that tracks it (Overview roll-up, zh-cn detail, de-de detail), the
# handoff timestamp advances, and the zh-cn/de-de "Priority" moves from
# "ht" to "mt". The Status column (and its Overview roll-up columns) also
# grow wider to fit the new text.

$wb = $excel.ActiveWorkbook

# ColumnWidth is specified in "characters"; Excel internally stores the
# resulting width snapped to a coarse pixel grid, so we pick the input
# that snaps closest to the target stored width (17.2159881591797).
$targetColWidth = 16.333333333333332

# --- Overview sheet: roll-up columns for the 631a... row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 08:14:53"
$wsOverview.Columns.Item(5).ColumnWidth = $targetColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColWidth

# --- zh-cn sheet: detail row for 631a... (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-19 08:14:48"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColWidth

# --- de-de sheet: detail row for 631a... (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-19 08:14:53"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColWidth
